$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 85.48
$ws.Range("E2").Value = "`n`nReasoning: The candidate has a strong skillset related to the job description, having worked with the MERN stack, ReactJS, NodeJS, ExpressJS, Socket.IO, WebRTC, HTML, CSS, and JS. The projects they have worked on demonstrate their ability to create and design web/mobile applications, manage data, and develop features with scalability. The score of 85.48 reflects their aptitude for the job."

# --- Row 3 updates (content that used to live in row 4, rephrased) ---
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 74.48
$ws.Range("E3").Value = "`n`nReasoning: The candidate has demonstrated excellent technical skills and knowledge of the MERN stack, Django Rest Framework, Pytorch, Tensorflow, Keras, and Sklearn. He has also shown experience in developing web/mobile applications, library management systems, user interfaces, and machine learning models. With these skills, he is well suited for the job and has been given a score of 74.48."

# --- Remove old row 4 entirely, shrinking the used range to A1:E3 ---
$ws.Rows("4:4").Delete()
